$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices / volume changes) and the row 37/38 swap
# Price cells that look like plain numbers must be forced to remain text (matching the
# original inlineStr cell type) by temporarily applying a Text number format, then
# reverting the cell style back to Normal so no stray style index is left on the cell.
$ws.Range('D2').Value = '36.029.58'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '1.914.49'
$ws.Range('E3').Value = '  -4.70%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.601'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.92%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.33'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -11.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.360'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '54.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0819'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.82%  '
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('D13').Value = '2.202.28'
$ws.Range('E13').Value = '  -4.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.802'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.01%  '
$ws.Range('D18').Value = '1.918.02'
$ws.Range('E18').Value = '  -4.54%  '
$ws.Range('D19').Value = '35.991.65'
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.01%  '
$ws.Range('D21').Value = '0.0₃0849'
$ws.Range('E21').Value = '  -2.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.52%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.118'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -13.90%  '
$ws.Range('E31').Value = '  -4.03%  '
$ws.Range('E32').Value = '  -5.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.59'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0617'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.70%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.90'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.86%  '
$ws.Range('E39').Value = '  -10.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.84'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0950'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.08%  '
$ws.Range('E42').Value = '  -3.41%  '
$ws.Range('E43').Value = '  -9.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0205'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.01%  '
$ws.Range('D46').Value = '1.326.83'
$ws.Range('E46').Value = '  -2.47%  '
$ws.Range('E47').Value = '  -9.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.13'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.03%  '
